$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "327.98"
    "E2" = "-1.23%"
    "G2" = "22"
    "D3" = "44.39"
    "E3" = "-0.77%"
    "G3" = "22"
    "D4" = "5.323"
    "E4" = "-3.92%"
    "G4" = "22"
    "D5" = "0.08369"
    "E5" = "1.86%"
    "G5" = "22"
    "D6" = "1.930"
    "E6" = "-5.96%"
    "G6" = "22"
    "E7" = "-0.69%"
    "G7" = "22"
    "E8" = "-3.16%"
    "G8" = "22"
    "D9" = "0.1133"
    "E9" = "-0.03%"
    "G9" = "22"
    "D10" = "0.1905"
    "E10" = "-0.41%"
    "G10" = "22"
    "D11" = "0.09650"
    "E11" = "-3.50%"
    "G11" = "22"
    "D12" = "0.04601"
    "E12" = "-1.83%"
    "G12" = "22"
    "E13" = "0.30%"
    "G13" = "22"
    "D14" = "0.001289"
    "E14" = "1.36%"
    "G14" = "22"
    "D15" = "0.005767"
    "E15" = "-2.03%"
    "G15" = "22"
    "E16" = "0.70%"
    "G16" = "22"
    "D17" = "4.402"
    "E17" = "-0.61%"
    "G17" = "22"
    "D18" = "0.3358"
    "G18" = "22"
    "D19" = "8.513"
    "E19" = "-17.59%"
    "G19" = "22"
    "D20" = "0.1389"
    "E20" = "1.00%"
    "G20" = "22"
    "D21" = "0.2576"
    "E21" = "3.46%"
    "G21" = "22"
    "D22" = "0.04153"
    "E22" = "1.21%"
    "G22" = "22"
    "E23" = "-5.35%"
    "G23" = "22"
    "D24" = "0.004407"
    "E24" = "-0.32%"
    "G24" = "22"
    "D25" = "0.0001299"
    "E25" = "1.54%"
    "G25" = "22"
    "D26" = "0.0002977"
    "E26" = "-20.37%"
    "G26" = "22"
    "G27" = "22"
    "G28" = "22"
    "G29" = "22"
    "G30" = "22"
    "G31" = "22"
    "G32" = "22"
    "G33" = "22"
    "G34" = "22"
    "G35" = "22"
    "G36" = "22"
    "G37" = "22"
    "D38" = "0.02714"
    "E38" = "-2.14%"
    "G38" = "22"
    "D39" = "0.05596"
    "E39" = "-2.58%"
    "G39" = "22"
    "D40" = "0.007852"
    "E40" = "2.81%"
    "G40" = "22"
    "D41" = "0.1414"
    "E41" = "-0.80%"
    "G41" = "22"
    "D42" = "0.007302"
    "E42" = "-3.49%"
    "G42" = "22"
    "D43" = "0.002116"
    "E43" = "7.28%"
    "G43" = "22"
    "D44" = "0.008687"
    "E44" = "4.56%"
    "G44" = "22"
    "D45" = "0.3514"
    "G45" = "22"
    "D46" = "0.00006900"
    "E46" = "-1.83%"
    "G46" = "22"
    "E47" = "-0.02%"
    "G47" = "22"
    "D48" = "0.003490"
    "E48" = "-1.89%"
    "G48" = "22"
    "D49" = "0.003528"
    "E49" = "40.05%"
    "G49" = "22"
    "E50" = "-0.02%"
    "G50" = "22"
    "E51" = "-0.02%"
    "G51" = "22"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
